$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" (column G) values
$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 1
    9  = 1
    10 = 3
    11 = 0
    12 = 1
    13 = 0
    14 = 2
    15 = 2
    16 = 3
    17 = 0
    18 = 2
    19 = 3
    21 = 1
    22 = 1
    23 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
